$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.935.83"
$ws.Range("E2").Value = "  +6.90%  "
$ws.Range("D3").Value = "3.014.18"
$ws.Range("E3").Value = "  +4.05%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'585.46"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").Value = "'155.78"
$ws.Range("E6").Value = "  +8.92%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.008.91"
$ws.Range("E8").Value = "  +3.96%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("D10").Value = "'6.94"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "'0.155"
$ws.Range("E11").Value = "  +6.37%  "
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("E13").Value = "  +8.28%  "
$ws.Range("D14").Value = "'34.56"
$ws.Range("E14").Value = "  +8.81%  "
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "65.910.76"
$ws.Range("E16").Value = "  +6.91%  "
$ws.Range("D17").Value = "3.513.92"
$ws.Range("E17").Value = "  +4.04%  "
$ws.Range("D18").Value = "'6.96"
$ws.Range("E18").Value = "  +6.75%  "
$ws.Range("D19").Value = "3.014.00"
$ws.Range("E19").Value = "  +3.99%  "
$ws.Range("D20").Value = "'463.89"
$ws.Range("E20").Value = "  +7.74%  "
$ws.Range("E21").Value = "  +6.21%  "
$ws.Range("D22").Value = "'0.684"
$ws.Range("E22").Value = "  +4.82%  "
$ws.Range("D24").Value = "'81.96"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "'2.26"
$ws.Range("E25").Value = "  +12.40%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'12.52"
$ws.Range("E26").Value = "  +5.59%  "
$ws.Range("D27").Value = "'10.70"
$ws.Range("E27").Value = "  +7.86%  "
$ws.Range("D29").Value = "'7.98"
$ws.Range("E29").Value = "  +13.84%  "
$ws.Range("D30").Value = "'2.38"
$ws.Range("E30").Value = "  +16.85%  "
$ws.Range("D31").Value = "'0.0000105"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("E32").Value = "  +5.23%  "
$ws.Range("D33").Value = "'0.112"
$ws.Range("E33").Value = "  +5.50%  "
$ws.Range("D34").Value = "'27.03"
$ws.Range("E34").Value = "  +5.90%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("E37").Value = "  +8.55%  "
$ws.Range("D38").Value = "'2.19"
$ws.Range("E38").Value = "  +13.54%  "
$ws.Range("E39").Value = "  +8.41%  "
$ws.Range("D40").Value = "'49.28"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.123"
$ws.Range("E41").Value = "  +8.65%  "
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").Value = "'44.28"
$ws.Range("E42").Value = "  +11.45%  "
$ws.Range("D43").Value = "'0.303"
$ws.Range("E43").Value = "  +13.85%  "
$ws.Range("D44").Value = "'8.45"
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("D45").Value = "'397.73"
$ws.Range("E45").Value = "  +15.96%  "
$ws.Range("D46").Value = "2.798.55"
$ws.Range("E46").Value = "  +4.28%  "
$ws.Range("D47").Value = "'0.0355"
$ws.Range("E47").Value = "  +5.74%  "
$ws.Range("D48").Value = "'133.84"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "'23.57"
$ws.Range("E50").Value = "  +9.55%  "
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  +3.65%  "
